$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 493.1579
$ws.Range("I6").Value = 23.846153
$ws.Range("J6").Value = 1510
$ws.Range("K6").Value = 71.538459
$ws.Range("L6").Value = 4530
$ws.Range("M6").Value = 40.461541
$ws.Range("N6").Value = -4754

$ws.Range("H20").Value = 1100
$ws.Range("I20").Value = 1100
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1100
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -870
$ws.Range("N20").ClearContents()

$ws.Range("H35").Value = 1100
$ws.Range("I35").Value = 1100
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1100
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -721
$ws.Range("N35").ClearContents()

$ws.Range("H121").Value = 878.75
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 878.75
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 2636.25
$ws.Range("N121").Value = -6130.25

$ws.Range("H132").Value = 11426.571
$ws.Range("I132").Value = 14303.9375
$ws.Range("J132").Value = 2219
$ws.Range("K132").Value = 42911.8125
$ws.Range("L132").Value = 6657
$ws.Range("M132").Value = -40381.8125
$ws.Range("N132").Value = -11717

$ws.Range("H141").Value = 2626.6
$ws.Range("I141").Value = 2008.25
$ws.Range("J141").Value = 5100
$ws.Range("K141").Value = 6024.75
$ws.Range("L141").Value = 15300
$ws.Range("M141").Value = -844.75
$ws.Range("N141").Value = -25660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4913.4443
$ws.Range("I32").Value = 4913.4443
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4913.4443
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4626.4443

$ws.Range("H33").Value = 11000.5
$ws.Range("I33").Value = 2001
$ws.Range("J33").Value = 20000
$ws.Range("K33").Value = 2001
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = -1672
$ws.Range("N33").Value = -20658

$ws.Range("H97").Value = 1405.8572
$ws.Range("I97").Value = 934.55554
$ws.Range("J97").Value = 2254.2
$ws.Range("K97").Value = 934.55554
$ws.Range("L97").Value = 2254.2
$ws.Range("M97").Value = -438.55554
$ws.Range("N97").Value = -3246.2

$ws.Range("H132").Value = 4671.143
$ws.Range("I132").Value = 1924.75
$ws.Range("J132").Value = 8333
$ws.Range("K132").Value = 5774.25
$ws.Range("L132").Value = 24999
$ws.Range("M132").Value = -3244.25
$ws.Range("N132").Value = -30059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 18498.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 18498.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 18498.5
$ws.Range("N88").Value = -19310.5

$ws.Range("H91").Value = 18498.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 18498.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 18498.5
$ws.Range("N91").Value = -21306.5

$ws.Range("H107").Value = 4543.4736
$ws.Range("I107").Value = 1393.2727
$ws.Range("J107").Value = 8875
$ws.Range("K107").Value = 1393.2727
$ws.Range("L107").Value = 8875
$ws.Range("M107").Value = 526.7273
$ws.Range("N107").Value = -12715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1362.0834
$ws.Range("I16").Value = 1326.909
$ws.Range("J16").Value = 1749
$ws.Range("K16").Value = 1326.909
$ws.Range("L16").Value = 1749
$ws.Range("M16").Value = -1039.909
$ws.Range("N16").Value = -2323

$ws.Range("H22").Value = 1214.6
$ws.Range("I22").Value = 268.25
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 268.25
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = 81.75
$ws.Range("N22").Value = -5700

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

$ws.Range("H58").Value = 2795.077
$ws.Range("I58").Value = 1599.875
$ws.Range("J58").Value = 4707.4
$ws.Range("K58").Value = 1599.875
$ws.Range("L58").Value = 4707.4
$ws.Range("M58").Value = -1396.875
$ws.Range("N58").Value = -5113.4

$ws.Range("H95").Value = 17507.273
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 17507.273
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 17507.273
$ws.Range("N95").Value = -22999.273

$ws.Range("H113").Value = 1362.0834
$ws.Range("I113").Value = 1326.909
$ws.Range("J113").Value = 1749
$ws.Range("K113").Value = 1326.909
$ws.Range("L113").Value = 1749
$ws.Range("M113").Value = 843.0909999999999
$ws.Range("N113").Value = -6089

$ws.Range("H134").Value = 1519.75
$ws.Range("I134").Value = 1519.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4559.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2024.25

$ws.Range("H136").Value = 2795.077
$ws.Range("I136").Value = 1599.875
$ws.Range("J136").Value = 4707.4
$ws.Range("K136").Value = 4799.625
$ws.Range("L136").Value = 14122.2
$ws.Range("M136").Value = -2249.625
$ws.Range("N136").Value = -19222.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1106.3077
$ws.Range("I5").Value = 1041.8572
$ws.Range("J5").Value = 1181.5
$ws.Range("K5").Value = 3125.5716
$ws.Range("L5").Value = 3544.5
$ws.Range("M5").Value = -3013.5716
$ws.Range("N5").Value = -3768.5

$ws.Range("H86").Value = 208.16667
$ws.Range("I86").Value = 198.6
$ws.Range("J86").Value = 256
$ws.Range("K86").Value = 595.8
$ws.Range("L86").Value = 768
$ws.Range("M86").Value = 590.2
$ws.Range("N86").Value = -3140

$ws.Range("H89").Value = 208.16667
$ws.Range("I89").Value = 198.6
$ws.Range("J89").Value = 256
$ws.Range("K89").Value = 1787.4
$ws.Range("L89").Value = 2304
$ws.Range("M89").Value = 4140.6
$ws.Range("N89").Value = -14160

$ws.Range("H122").Value = 620.3333
$ws.Range("I122").Value = 508.33334
$ws.Range("J122").Value = 732.3333
$ws.Range("K122").Value = 4575.00006
$ws.Range("L122").Value = 6590.9997
$ws.Range("M122").Value = -2125.00006
$ws.Range("N122").Value = -11490.9997

$ws.Range("H132").Value = 2350
$ws.Range("I132").Value = 2101
$ws.Range("J132").Value = 2516
$ws.Range("K132").Value = 18909
$ws.Range("L132").Value = 22644
$ws.Range("M132").Value = -16379
$ws.Range("N132").Value = -27704

$ws.Range("H135").Value = 1106.3077
$ws.Range("I135").Value = 1041.8572
$ws.Range("J135").Value = 1181.5
$ws.Range("K135").Value = 9376.7148
$ws.Range("L135").Value = 10633.5
$ws.Range("M135").Value = -6841.7148
$ws.Range("N135").Value = -15703.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11858974
$ws.Range("I11").Value = 9931034
$ws.Range("J11").Value = 17450000
$ws.Range("K11").Value = 9931034
$ws.Range("L11").Value = 17450000
$ws.Range("M11").Value = -9930895
$ws.Range("N11").Value = -17450278

$ws.Range("H42").Value = 79750
$ws.Range("I42").Value = 75000
$ws.Range("J42").Value = 84500
$ws.Range("K42").Value = 75000
$ws.Range("L42").Value = 84500
$ws.Range("M42").Value = -74515
$ws.Range("N42").Value = -85470

$ws.Range("H97").Value = 531.4211
$ws.Range("I97").Value = 556.94446
$ws.Range("J97").Value = 72
$ws.Range("K97").Value = 556.94446
$ws.Range("L97").Value = 72
$ws.Range("M97").Value = -60.94446000000005
$ws.Range("N97").Value = -1064

$ws.Range("H115").Value = 79750
$ws.Range("I115").Value = 75000
$ws.Range("J115").Value = 84500
$ws.Range("K115").Value = 75000
$ws.Range("L115").Value = 84500
$ws.Range("M115").Value = -73825
$ws.Range("N115").Value = -86850

$ws.Range("H135").Value = 240000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 240000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 240000
$ws.Range("N135").Value = -250140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 4992
$ws.Range("I10").Value = 4999
$ws.Range("J10").Value = 4985
$ws.Range("K10").Value = 4999
$ws.Range("L10").Value = 4985
$ws.Range("M10").Value = -4859
$ws.Range("N10").Value = -5265

$ws.Range("H16").Value = 999.3333
$ws.Range("I16").Value = 999.3333
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 999.3333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -829.3333

$ws.Range("H18").Value = 1166.6666
$ws.Range("I18").Value = 1100
$ws.Range("J18").Value = 1200
$ws.Range("K18").Value = 1100
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = -928
$ws.Range("N18").Value = -1544

$ws.Range("H40").Value = 7069.5
$ws.Range("I40").Value = 6223.2
$ws.Range("J40").Value = 9185.25
$ws.Range("K40").Value = 6223.2
$ws.Range("L40").Value = 9185.25
$ws.Range("M40").Value = -6087.2
$ws.Range("N40").Value = -9457.25

$ws.Range("H61").Value = 6988.0835
$ws.Range("I61").Value = 5977.8335
$ws.Range("J61").Value = 7998.3335
$ws.Range("K61").Value = 5977.8335
$ws.Range("L61").Value = 7998.3335
$ws.Range("M61").Value = -5775.8335
$ws.Range("N61").Value = -8402.333500000001

$ws.Range("H68").Value = 10100.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 10100.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 10100.5
$ws.Range("N68").Value = -11598.5

$ws.Range("H71").Value = 10100.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 10100.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 50502.5
$ws.Range("N71").Value = -57990.5

$ws.Range("H113").Value = 6988.0835
$ws.Range("I113").Value = 5977.8335
$ws.Range("J113").Value = 7998.3335
$ws.Range("K113").Value = 5977.8335
$ws.Range("L113").Value = 7998.3335
$ws.Range("M113").Value = -3807.8335
$ws.Range("N113").Value = -12338.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("N34").ClearContents()

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").ClearContents()
